# Add three more days (columns JS, JT, JU) of mobility data to the "mobility" sheet,
# mirroring the formatting of the last existing date column (JR).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-number-format from JR1 onto the three new header cells (JS1:JU1)
# so the new style reuses the existing cellXf (numFmtId 14) instead of creating a new one.
$ws.Range("JR1").Copy()
$ws.Range("JS1:JU1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New date header values (serial dates continuing 44114 -> 44115, 44116, 44117)
$ws.Range("JS1").Value2 = 44115
$ws.Range("JT1").Value2 = 44116
$ws.Range("JU1").Value2 = 44117

# Row 2 values
$ws.Range("JS2").Value2 = 48.49
$ws.Range("JT2").Value2 = 42.43
$ws.Range("JU2").Value2 = 47.13

# Row 3 values
$ws.Range("JS3").Value2 = 35.44
$ws.Range("JT3").Value2 = 37.840000000000003
$ws.Range("JU3").Value2 = 36.409999999999997

# Row 4 values
$ws.Range("JS4").Value2 = 46.18
$ws.Range("JT4").Value2 = 57.06
$ws.Range("JU4").Value2 = 57.56

# Row 5 values
$ws.Range("JS5").Value2 = 50.65
$ws.Range("JT5").Value2 = 61.31
$ws.Range("JU5").Value2 = 63.05

# Mirror the view state change recorded in the diff (scroll position / active selection).
[void]$ws.Range("JJ1").Select()
[void]$ws.Range("JZ22").Select()
